# Updates the cryptos price list (Coin/Link/Price/Volume(1h)) to the
# latest scraped snapshot. Columns D (Price) and E (Volume 1h) are
# stored as text in this sheet (not numbers), so any value that Excel
# would otherwise auto-convert to a number on assignment is written
# with a leading apostrophe to force text, matching the original
# "t=inlineStr" layout of the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.146.73"
$ws.Range("E2").Value = "  +0.24%  "

$ws.Range("D3").Value = "4.037.28"
$ws.Range("E3").Value = "  -0.18%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "'538.70"
$ws.Range("E5").Value = "  -0.62%  "

$ws.Range("D6").Value = "'149.03"
$ws.Range("E6").Value = "  -1.14%  "

$ws.Range("D7").Value = "4.030.50"
$ws.Range("E7").Value = "  -0.20%  "

$ws.Range("D8").Value = "'0.695"
$ws.Range("E8").Value = "  +0.12%  "

$ws.Range("E9").Value = "  -0.01%  "

$ws.Range("D10").Value = "'0.750"
$ws.Range("E10").Value = "  -1.76%  "

$ws.Range("D11").Value = "'0.173"
$ws.Range("E11").Value = "  -1.31%  "

$ws.Range("D12").Value = "'53.02"
$ws.Range("E12").Value = "  +10.32%  "

$ws.Range("E13").Value = "  -0.77%  "

$ws.Range("D14").Value = "'10.86"
$ws.Range("E14").Value = "  -0.78%  "

$ws.Range("D15").Value = "4.684.14"
$ws.Range("E15").Value = "  +0.03%  "

$ws.Range("D16").Value = "4.054.21"
$ws.Range("E16").Value = "  +0.54%  "

$ws.Range("D17").Value = "'14.27"
$ws.Range("E17").Value = "  -1.18%  "

$ws.Range("D18").Value = "'20.64"
$ws.Range("E18").Value = "  -0.27%  "

$ws.Range("E19").Value = "  -0.85%  "

$ws.Range("E20").Value = "  -1.09%  "

$ws.Range("D21").Value = "72.096.92"
$ws.Range("E21").Value = "  +0.50%  "

$ws.Range("D22").Value = "'438.29"
$ws.Range("E22").Value = "  +0.56%  "

$ws.Range("D23").Value = "'98.01"
$ws.Range("E23").Value = "  -1.74%  "

$ws.Range("D24").Value = "'3.49"
$ws.Range("E24").Value = "  -3.85%  "

$ws.Range("D25").Value = "'4.27"
$ws.Range("E25").Value = "  -1.06%  "

$ws.Range("D26").Value = "'14.59"
$ws.Range("E26").Value = "  -1.06%  "

$ws.Range("D27").Value = "'4.41"
$ws.Range("E27").Value = "  +26.72%  "

$ws.Range("E28").Value = "  -0.74%  "

$ws.Range("D29").Value = "'10.67"
$ws.Range("E29").Value = "  -2.21%  "

$ws.Range("D30").Value = "'5.95"
$ws.Range("E30").Value = "  +1.77%  "

$ws.Range("D31").Value = "'37.08"
$ws.Range("E31").Value = "  -0.56%  "

$ws.Range("D32").Value = "'8.37"
$ws.Range("E32").Value = "  +22.96%  "

$ws.Range("D33").Value = "'0.134"
$ws.Range("E33").Value = "  +1.86%  "

$ws.Range("D34").Value = "'13.55"
$ws.Range("E34").Value = "  -0.48%  "

$ws.Range("D35").Value = "'49.49"
$ws.Range("E35").Value = "  +15.89%  "

$ws.Range("D36").Value = "'682.22"
$ws.Range("E36").Value = "  -0.32%  "

$ws.Range("D37").Value = "'66.81"
$ws.Range("E37").Value = "  -0.40%  "

$ws.Range("D38").Value = "'0.458"
$ws.Range("E38").Value = "  +4.59%  "

$ws.Range("D39").Value = "0.0₃0899"
$ws.Range("E39").Value = "  +5.57%  "

$ws.Range("D40").Value = "'3.43"
$ws.Range("E40").Value = "  +6.05%  "

$ws.Range("B41").Value = "THORChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D41").Value = "'11.36"
$ws.Range("E41").Value = "  +18.43%  "

$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "'0.148"
$ws.Range("E42").Value = "  -6.16%  "

$ws.Range("B43").Value = "ThetaToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D43").Value = "'3.39"
$ws.Range("E43").Value = "  -2.20%  "

$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "  -0.03%  "

$ws.Range("D45").Value = "'0.0493"
$ws.Range("E45").Value = "  -0.89%  "

$ws.Range("E46").Value = "  +0.14%  "

$ws.Range("D47").Value = "'0.150"
$ws.Range("E47").Value = "  -1.64%  "

$ws.Range("D48").Value = "'2.64"
$ws.Range("E48").Value = "  -3.64%  "

$ws.Range("D49").Value = "'3.12"
$ws.Range("E49").Value = "  +2.29%  "

$ws.Range("D50").Value = "'3.32"
$ws.Range("E50").Value = "  -3.39%  "

$ws.Range("D51").Value = "'0.000278"
$ws.Range("E51").Value = "  +1.35%  "
